# This workbook is a weekly price log. A new week's worth of observations
# (two rows: "Primera" and "Segunda" quality) is inserted at the top of the
# data block that starts at row 408, pushing the existing rows down by two
# positions (408-466 -> 410-468).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 408; this shifts old rows 408:466 down
# to 410:468 and carries their formatting (e.g. the date style on column D)
# along with them.
$ws.Rows("408:409").Insert()

# Fixed values shared by every row in this subset.
$market      = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$catId       = 100112009
$categoria   = "Acelga"
$variedad    = "Sin especificar"
$unidad      = "$/atado 0,5 a 1 kilo"
$origen      = "Provincia de Diguillín"
$clasif      = "Hortaliza"

# New row 408: "Primera" quality observation for the new week (45077).
$r = 408
$ws.Cells.Item($r, 1).Value = 7
$ws.Cells.Item($r, 2).Value = $market
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45077
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 160
$ws.Cells.Item($r, 11).Value = 600
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 650
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 650
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = $clasif

# New row 409: "Segunda" quality observation for the same new week.
$r = 409
$ws.Cells.Item($r, 1).Value = 7
$ws.Cells.Item($r, 2).Value = $market
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45077
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 150
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = $clasif
